$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '63.177.13'
$ws.Range('E2').Value = '  -0.41%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '3.049.99'
$ws.Range('E3').Value = '  -0.69%  '
$ws.Range('E4').Value = '  -0.07%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '588.03'
$ws.Range('E5').Value = '  -0.30%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '151.42'
$ws.Range('E6').Value = '  -0.88%  '
$ws.Range('E7').Value = '  -0.03%  '
$ws.Range('E8').Value = '  -1.44%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '3.050.06'
$ws.Range('E9').Value = '  -0.68%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.154'
$ws.Range('E10').Value = '  -1.31%  '
$ws.Range('E11').Value = '  -0.60%  '
$ws.Range('E12').Value = '  -2.47%  '
$ws.Range('E13').Value = '  -2.44%  '
$ws.Range('E14').Value = '  -2.10%  '
$ws.Range('E15').Value = '  +1.67%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '3.552.97'
$ws.Range('E16').Value = '  -0.77%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '7.17'
$ws.Range('E17').Value = '  -0.57%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '63.202.48'
$ws.Range('E18').Value = '  -0.38%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '3.048.43'
$ws.Range('E19').Value = '  -0.83%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '477.82'
$ws.Range('E20').Value = '  +0.64%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '14.28'
$ws.Range('E21').Value = '  -2.62%  '
$ws.Range('E22').Value = '  -1.48%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '7.52'
$ws.Range('E23').Value = '  -0.32%  '
$ws.Range('E24').Value = '  +1.64%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '82.13'
$ws.Range('E25').Value = '  +0.98%  '
$ws.Range('E26').Value = '  -2.35%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '10.73'
$ws.Range('E27').Value = '  +7.24%  '
$ws.Range('E28').Value = '  +0.30%  '
$ws.Range('E29').Value = '  +0.55%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '2.68'
$ws.Range('E30').Value = '  -0.13%  '
$ws.Range('E31').Value = '  -0.20%  '
$ws.Range('E32').Value = '  +0.14%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '27.70'
$ws.Range('E33').Value = '  +1.48%  '
$ws.Range('E34').Value = '  -2.50%  '
$ws.Range('E35').Value = '  +1.00%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.0₃0818'
$ws.Range('E36').Value = '  -4.05%  '
$ws.Range('B37').Value = 'Stacks'
$ws.Range('C37').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '2.22'
$ws.Range('E37').Value = '  +0.56%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '5.90'
$ws.Range('E38').Value = '  -3.53%  '
$ws.Range('B39').Value = 'dogwifhat'
$ws.Range('C39').Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '3.24'
$ws.Range('E39').Value = '  -3.46%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '9.28'
$ws.Range('E40').Value = '  -0.30%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '50.52'
$ws.Range('E41').Value = '  +0.13%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '436.12'
$ws.Range('E42').Value = '  -2.65%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.288'
$ws.Range('E43').Value = '  +0.76%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.114'
$ws.Range('E44').Value = '  +2.70%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.0363'
$ws.Range('E45').Value = '  -0.20%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '2.829.71'
$ws.Range('E46').Value = '  +0.85%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '38.30'
$ws.Range('E47').Value = '  -4.78%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '129.50'
$ws.Range('E48').Value = '  -1.28%  '
$ws.Range('E49').Value = '  +0.00%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '25.12'
$ws.Range('E50').Value = '  -0.40%  '
$ws.Range('B51').Value = 'Stellar'
$ws.Range('C51').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.109'
$ws.Range('E51').Value = '  -1.72%  '
